# Apply updated cryptocurrency price/volume data per diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.702.37"
$ws.Range("E2").Value = "  -4.67%  "

# Row 3
$ws.Range("D3").Value = "3.465.42"
$ws.Range("E3").Value = "  -6.19%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.44"
$ws.Range("E5").Value = "  -7.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.59"
$ws.Range("E6").Value = "  -9.01%  "

# Row 7
$ws.Range("D7").Value = "3.461.94"
$ws.Range("E7").Value = "  -6.24%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("E9").Value = "  -4.55%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  -6.39%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.89"
$ws.Range("E11").Value = "  -4.52%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.419"
$ws.Range("E12").Value = "  -5.64%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000214"
$ws.Range("E13").Value = "  -7.94%  "

# Row 14
$ws.Range("D14").Value = "4.036.45"
$ws.Range("E14").Value = "  -6.54%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.05"
$ws.Range("E15").Value = "  -5.40%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.459.31"
$ws.Range("E16").Value = "  -6.33%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.564.72"
$ws.Range("E17").Value = "  -4.87%  "

# Row 18
$ws.Range("E18").Value = "  -0.80%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("E19").Value = "  -1.92%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.90"
$ws.Range("E20").Value = "  -6.91%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "441.41"
$ws.Range("E21").Value = "  -6.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.97"
$ws.Range("E22").Value = "  -13.88%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.622"
$ws.Range("E23").Value = "  -5.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.78"
$ws.Range("E24").Value = "  -4.23%  "

# Row 25
$ws.Range("E25").Value = "  +0.11%  "

# Row 26
$ws.Range("D26").Value = "3.590.25"
$ws.Range("E26").Value = "  -6.52%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000124"
$ws.Range("E27").Value = "  -2.77%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  -8.68%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.26"
$ws.Range("E29").Value = "  -9.74%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.51"
$ws.Range("E30").Value = "  -5.73%  "

# Row 31
$ws.Range("E31").Value = "  -0.09%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.54"
$ws.Range("E32").Value = "  -10.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.160"
$ws.Range("E33").Value = "  -5.51%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.52"
$ws.Range("E34").Value = "  -4.74%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.11"
$ws.Range("E35").Value = "  -6.75%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.84"
$ws.Range("E36").Value = "  -8.51%  "

# Row 37
$ws.Range("D37").Value = "3.437.89"
$ws.Range("E37").Value = "  -6.85%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.91"
$ws.Range("E38").Value = "  -6.61%  "

# Row 39
$ws.Range("E39").Value = "  +0.11%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  -0.28%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "172.64"
$ws.Range("E41").Value = "  -4.38%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.13"
$ws.Range("E42").Value = "  -4.58%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0863"
$ws.Range("E43").Value = "  -4.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.41"
$ws.Range("E44").Value = "  -8.30%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.875"
$ws.Range("E45").Value = "  -6.20%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.45"
$ws.Range("E46").Value = "  -2.46%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.22"
$ws.Range("E47").Value = "  -3.11%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.82"
$ws.Range("E48").Value = "  -12.54%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.54"
$ws.Range("E49").Value = "  -4.27%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.48"
$ws.Range("E50").Value = "  -13.45%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -5.56%  "
